$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used row in the sheet so we cover the whole table.
$lastRow = $ws.UsedRange.Rows.Count

# Column G is the "Recorded By" column, containing comma-separated
# recorder names/emails (e.g. "System, dnasr281@gmail.com").
# Reverse the order of the items in each cell's list.
$col = 7

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $col)
    $val = $cell.Value2

    if ($val -ne $null -and $val -is [string] -and $val.Contains(",")) {
        $parts = $val -split ", "
        if ($parts.Count -gt 1) {
            $reversed = $parts[($parts.Count - 1)..0]
            $newVal = $reversed -join ", "
            $cell.Value2 = $newVal
        }
    }
}
